$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6485
$ws.Range("B2").Value = "Eloah Santos"
$ws.Range("C2").Value = "Operações"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45086
$ws.Range("G2").Value = 6974.91

$ws.Range("A3").Value = 79361
$ws.Range("B3").Value = "Samuel Gomes"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45097
$ws.Range("G3").Value = 8809.11

$ws.Range("A4").Value = 70045
$ws.Range("B4").Value = "Danilo Barros"
$ws.Range("C4").Value = "P&D"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45091
$ws.Range("G4").Value = 9920.58

$ws.Range("A5").Value = 75867
$ws.Range("B5").Value = "Isadora Ferreira"
$ws.Range("C5").Value = "Recursos Humanos"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45080
$ws.Range("G5").Value = 11851.92

$ws.Range("A6").Value = 16037
$ws.Range("B6").Value = "Vitor Gabriel Monteiro"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 10133.61

$ws.Range("A7").Value = 22265
$ws.Range("B7").Value = "Davi Luiz Melo"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45091
$ws.Range("G7").Value = 9473.27

$ws.Range("A8").Value = 61272
$ws.Range("B8").Value = "Emanuel Dias"
$ws.Range("C8").Value = "Operações"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45101
$ws.Range("G8").Value = 6428.78

$ws.Range("A9").Value = 20069
$ws.Range("B9").Value = "Cauê Rezende"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45086
$ws.Range("G9").Value = 4189.75

$ws.Range("A10").Value = 23058
$ws.Range("B10").Value = "Helena Cunha"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45082
$ws.Range("G10").Value = 12038.77

$ws.Range("A11").Value = 64271
$ws.Range("B11").Value = "Dra. Melissa da Costa"
$ws.Range("C11").Value = "Recursos Humanos"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45105
$ws.Range("G11").Value = 6150.01

Write-Host "Updated rows 2-11 of absenteeism data (ETL refresh)"
